# Update stock data in spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement Stock (column B) for the rows that were restocked/sold.
$ws.Range("B5").Value = 2
$ws.Range("B19").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("B22").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B34").Value = 1
$ws.Range("B35").Value = 2
$ws.Range("B48").Value = 1
$ws.Range("B51").Value = 0
$ws.Range("B59").Value = 1
$ws.Range("B63").Value = 4
$ws.Range("B65").Value = 1
$ws.Range("B70").Value = 2
$ws.Range("B86").Value = 1
$ws.Range("B240").Value = 0
$ws.Range("B244").Value = 0
$ws.Range("B247").Value = 0
$ws.Range("B249").Value = 0
$ws.Range("B279").Value = 3
$ws.Range("B294").Value = 0

# Update the active view/selection to match the scrolled position in the
# saved workbook.
$ws.Range("B249").Select()
$excel.ActiveWindow.ScrollRow = 241
